$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("per_game")
$ws.Range("G44").Value = 67
$ws.Range("H44").Value = 67
$ws.Range("L44").Value = 0.428
$ws.Range("O44").Value = 0.398
$ws.Range("R44").Value = 0.474
$ws.Range("S44").Value = 0.549
$ws.Range("U44").Value = 0.3
$ws.Range("AA44").Value = 0.6
$ws.Range("G55").Value = 67
$ws.Range("H55").Value = 67
$ws.Range("L55").Value = 0.428
$ws.Range("O55").Value = 0.398
$ws.Range("R55").Value = 0.474
$ws.Range("S55").Value = 0.549
$ws.Range("U55").Value = 0.3
$ws.Range("AA55").Value = 0.6
$ws.Range("G62").Value = 844
$ws.Range("H62").Value = 648
$ws.Range("I62").Value = 28.7
$ws.Range("K62").Value = 5.8
$ws.Range("K64").Value = 0.1000000000000005

$ws = $wb.Worksheets.Item("per_minute")
$ws.Range("G34").Value = 67
$ws.Range("H34").Value = 67
$ws.Range("I34").Value = 1726
$ws.Range("K34").Value = 4
$ws.Range("L34").Value = 0.428
$ws.Range("O34").Value = 0.398
$ws.Range("R34").Value = 0.474
$ws.Range("Z34").Value = 0.8
$ws.Range("AB34").Value = 0.8
$ws.Range("AC34").Value = 3.4
$ws.Range("AD34").Value = 4.8
$ws.Range("G45").Value = 67
$ws.Range("H45").Value = 67
$ws.Range("I45").Value = 1726
$ws.Range("K45").Value = 4
$ws.Range("L45").Value = 0.428
$ws.Range("O45").Value = 0.398
$ws.Range("R45").Value = 0.474
$ws.Range("Z45").Value = 0.8
$ws.Range("AB45").Value = 0.8
$ws.Range("AC45").Value = 3.4
$ws.Range("AD45").Value = 4.8
$ws.Range("G52").Value = 844
$ws.Range("H52").Value = 648
$ws.Range("I52").Value = 24255

$ws = $wb.Worksheets.Item("per_poss")
$ws.Range("G34").Value = 67
$ws.Range("H34").Value = 67
$ws.Range("I34").Value = 1726
$ws.Range("L34").Value = 0.428
$ws.Range("M34").Value = 1.3
$ws.Range("O34").Value = 0.398
$ws.Range("R34").Value = 0.474
$ws.Range("S34").Value = 0.5
$ws.Range("Z34").Value = 1.1
$ws.Range("AB34").Value = 1.1
$ws.Range("AF34").Value = 118
$ws.Range("AG34").Value = 115
$ws.Range("G45").Value = 67
$ws.Range("H45").Value = 67
$ws.Range("I45").Value = 1726
$ws.Range("L45").Value = 0.428
$ws.Range("M45").Value = 1.3
$ws.Range("O45").Value = 0.398
$ws.Range("R45").Value = 0.474
$ws.Range("S45").Value = 0.5
$ws.Range("Z45").Value = 1.1
$ws.Range("AB45").Value = 1.1
$ws.Range("AF45").Value = 118
$ws.Range("AG45").Value = 115
$ws.Range("G52").Value = 844
$ws.Range("H52").Value = 648
$ws.Range("I52").Value = 24255

$ws = $wb.Worksheets.Item("advanced")
$ws.Range("G34").Value = 67
$ws.Range("H34").Value = 1726
$ws.Range("I34").Value = 5.8
$ws.Range("J34").Value = 0.568
$ws.Range("K34").Value = 0.608
$ws.Range("L34").Value = 0.119
$ws.Range("M34").Value = 6.1
$ws.Range("N34").Value = 12.3
$ws.Range("Q34").Value = 1.1
$ws.Range("S34").Value = 16.4
$ws.Range("T34").Value = 6.3
$ws.Range("W34").Value = 1.7
$ws.Range("X34").Value = 2.5
$ws.Range("Y34").Value = 0.07
$ws.Range("AB34").Value = 0.4
$ws.Range("AC34").Value = -2.8
$ws.Range("G45").Value = 67
$ws.Range("H45").Value = 1726
$ws.Range("I45").Value = 5.8
$ws.Range("J45").Value = 0.568
$ws.Range("K45").Value = 0.608
$ws.Range("L45").Value = 0.119
$ws.Range("M45").Value = 6.1
$ws.Range("N45").Value = 12.3
$ws.Range("Q45").Value = 1.1
$ws.Range("S45").Value = 16.4
$ws.Range("T45").Value = 6.3
$ws.Range("W45").Value = 1.7
$ws.Range("X45").Value = 2.5
$ws.Range("Y45").Value = 0.07
$ws.Range("AB45").Value = 0.4
$ws.Range("AC45").Value = -2.8
$ws.Range("G52").Value = 844
$ws.Range("H52").Value = 24255
$ws.Range("S52").Value = 13
$ws.Range("W52").Value = 24.4
$ws.Range("X52").Value = 42.5
$ws.Range("AD52").Value = 6.8
$ws.Range("S54").Value = -2.800000000000001
$ws.Range("W54").Value = -20.7
$ws.Range("X54").Value = -35.5
$ws.Range("AD54").Value = -4.8
